$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the organization website URL (shared string previously "www.stat.kg")
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move/update the active selection to B10, matching the saved view state
$ws.Range("B10").Select()

# Reflect the updated window placement/size recorded by Excel on save
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 0
$win.Width = 1440
$win.Height = 591.75
